$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1428.8
$ws.Range("J12").Value = 2999.5
$ws.Range("L12").Value = 2999.5
$ws.Range("N12").Value = -3339.5
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H33").Value = 19231020
$ws.Range("I33").Value = 25000118
$ws.Range("J33").Value = 696
$ws.Range("K33").Value = 25000118
$ws.Range("L33").Value = 696
$ws.Range("M33").Value = -24999889
$ws.Range("N33").Value = -1154
$ws.Range("H62").Value = 2653.0908
$ws.Range("I62").Value = 2649.25
$ws.Range("K62").Value = 2649.25
$ws.Range("M62").Value = -2025.25
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H65").Value = 2653.0908
$ws.Range("I65").Value = 2649.25
$ws.Range("K65").Value = 13246.25
$ws.Range("M65").Value = -10126.25
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H76").Value = 17897.03
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 17897.03
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H82").Value = 1084.4286
$ws.Range("I82").Value = 1084.4286
$ws.Range("K82").Value = 3253.2858
$ws.Range("M82").Value = -2847.2858
$ws.Range("H85").Value = 1084.4286
$ws.Range("I85").Value = 1084.4286
$ws.Range("K85").Value = 3253.2858
$ws.Range("M85").Value = -1849.2858
$ws.Range("H97").Value = 1000
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992
$ws.Range("H99").Value = 191.33333
$ws.Range("I99").Value = 162
$ws.Range("J99").Value = 250
$ws.Range("K99").Value = 486
$ws.Range("L99").Value = 750
$ws.Range("M99").Value = 1012
$ws.Range("N99").Value = -3746
$ws.Range("H132").Value = 1791.2368
$ws.Range("I132").Value = 1451.2258
$ws.Range("K132").Value = 4353.6774
$ws.Range("M132").Value = -1823.6774
$ws.Range("H134").Value = 178015
$ws.Range("J134").Value = 178015
$ws.Range("L134").Value = 178015
$ws.Range("N134").Value = -188155
$ws.Range("H135").Value = 1949.125
$ws.Range("I135").Value = 999.2
$ws.Range("J135").Value = 3532.3333
$ws.Range("K135").Value = 8992.800000000001
$ws.Range("L135").Value = 31790.9997
$ws.Range("M135").Value = -6457.800000000001
$ws.Range("N135").Value = -36860.9997
$ws.Range("H137").Value = 3034.1875
$ws.Range("J137").Value = 3949.875
$ws.Range("L137").Value = 11849.625
$ws.Range("N137").Value = -16949.625
$ws.Range("H138").Value = 15626854
$ws.Range("I138").Value = 981.8837
$ws.Range("J138").Value = 47622684
$ws.Range("K138").Value = 2945.6511
$ws.Range("L138").Value = 142868052
$ws.Range("M138").Value = 2194.3489
$ws.Range("N138").Value = -142878332
$ws.Range("H141").Value = 2894.318
$ws.Range("I141").Value = 2989.3333
$ws.Range("K141").Value = 8967.999899999999
$ws.Range("M141").Value = -3787.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 500
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 500
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -350
$ws.Range("N30").ClearContents()
$ws.Range("H32").Value = 13260.088
$ws.Range("I32").Value = 3823.6538
$ws.Range("K32").Value = 3823.6538
$ws.Range("M32").Value = -3536.6538
$ws.Range("H46").Value = 9178.571
$ws.Range("I46").Value = 5949.5
$ws.Range("J46").Value = 10470.2
$ws.Range("K46").Value = 5949.5
$ws.Range("L46").Value = 10470.2
$ws.Range("M46").Value = -5630.5
$ws.Range("N46").Value = -11108.2
$ws.Range("H61").Value = 20844878
$ws.Range("I61").Value = 27788724
$ws.Range("K61").Value = 27788724
$ws.Range("M61").Value = -27788512
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H74").Value = 4072.7273
$ws.Range("I74").Value = 2707.6924
$ws.Range("K74").Value = 2707.6924
$ws.Range("M74").Value = -1833.6924
$ws.Range("H77").Value = 4072.7273
$ws.Range("I77").Value = 2707.6924
$ws.Range("K77").Value = 13538.462
$ws.Range("M77").Value = -9170.462
$ws.Range("H97").Value = 882.1667
$ws.Range("I97").Value = 946
$ws.Range("J97").Value = 639.6
$ws.Range("K97").Value = 946
$ws.Range("L97").Value = 639.6
$ws.Range("M97").Value = -450
$ws.Range("N97").Value = -1631.6
$ws.Range("H110").Value = 4844.1714
$ws.Range("I110").Value = 3903.5356
$ws.Range("K110").Value = 3903.5356
$ws.Range("M110").Value = -1858.5356
$ws.Range("H132").Value = 2729.5386
$ws.Range("I132").Value = 2698.6572
$ws.Range("J132").Value = 2999.75
$ws.Range("K132").Value = 8095.971600000001
$ws.Range("L132").Value = 8999.25
$ws.Range("M132").Value = -5565.971600000001
$ws.Range("N132").Value = -14059.25
$ws.Range("H136").Value = 20844878
$ws.Range("I136").Value = 27788724
$ws.Range("K136").Value = 83366172
$ws.Range("M136").Value = -83363622

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H86").Value = 2881.0908
$ws.Range("I86").Value = 2881.0908
$ws.Range("K86").Value = 2881.0908
$ws.Range("M86").Value = -1758.0908
$ws.Range("H89").Value = 2881.0908
$ws.Range("I89").Value = 2881.0908
$ws.Range("K89").Value = 14405.454
$ws.Range("M89").Value = -8789.454
$ws.Range("H94").Value = 2618.7778
$ws.Range("I94").Value = 2071.125
$ws.Range("K94").Value = 2071.125
$ws.Range("M94").Value = -1620.125
$ws.Range("H105").Value = 2309.75
$ws.Range("I105").Value = 2354
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2354
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -607
$ws.Range("N105").Value = -5494
$ws.Range("H130").Value = 40000
$ws.Range("I130").Value = 40000
$ws.Range("K130").Value = 40000
$ws.Range("M130").Value = -34980
$ws.Range("H134").Value = 4146.073
$ws.Range("J134").Value = 3607.5
$ws.Range("L134").Value = 10822.5
$ws.Range("N134").Value = -15892.5
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 35000
$ws.Range("I42").Value = 35000
$ws.Range("K42").Value = 35000
$ws.Range("M42").Value = -34407
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52372
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -161856
$ws.Range("H94").Value = 1764.7333
$ws.Range("J94").Value = 2048.4285
$ws.Range("L94").Value = 2048.4285
$ws.Range("N94").Value = -2950.4285
$ws.Range("H97").Value = 38000
$ws.Range("J97").Value = 38000
$ws.Range("L97").Value = 38000
$ws.Range("N97").Value = -39982
$ws.Range("H99").Value = 5607.1
$ws.Range("I99").Value = 5134
$ws.Range("K99").Value = 5134
$ws.Range("M99").Value = -3636
$ws.Range("H105").Value = 1821.1666
$ws.Range("I105").Value = 1395.5
$ws.Range("K105").Value = 1395.5
$ws.Range("M105").Value = 351.5
$ws.Range("H126").Value = 5607.1
$ws.Range("I126").Value = 5134
$ws.Range("K126").Value = 15402
$ws.Range("M126").Value = -12932
$ws.Range("H132").Value = 2769.4849
$ws.Range("I132").Value = 3021.1482
$ws.Range("J132").Value = 1637
$ws.Range("K132").Value = 9063.444600000001
$ws.Range("L132").Value = 4911
$ws.Range("M132").Value = -6533.444600000001
$ws.Range("N132").Value = -9971
$ws.Range("H134").Value = 3500.0212
$ws.Range("I134").Value = 2638.027
$ws.Range("J134").Value = 6689.4
$ws.Range("K134").Value = 7914.081
$ws.Range("L134").Value = 20068.2
$ws.Range("M134").Value = -5379.081
$ws.Range("N134").Value = -25138.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1047.25
$ws.Range("I8").Value = 1047.25
$ws.Range("K8").Value = 3141.75
$ws.Range("M8").Value = -3002.75
$ws.Range("H113").Value = 1041.25
$ws.Range("J113").Value = 1431
$ws.Range("L113").Value = 4293
$ws.Range("N113").Value = -8633
$ws.Range("H132").Value = 38463304
$ws.Range("J132").Value = 1997
$ws.Range("L132").Value = 17973
$ws.Range("N132").Value = -23033
$ws.Range("H140").Value = 454174.78
$ws.Range("I140").Value = 987.0526
$ws.Range("K140").Value = 2961.1578
$ws.Range("M140").Value = 2218.8422

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 9441.5
$ws.Range("J19").Value = 11976
$ws.Range("L19").Value = 11976
$ws.Range("N19").Value = -12552
$ws.Range("H121").Value = 60330.8
$ws.Range("J121").Value = 60330.8
$ws.Range("L121").Value = 60330.8
$ws.Range("N121").Value = -63824.8
$ws.Range("H122").Value = 2023.6154
$ws.Range("I122").Value = 2177.75
$ws.Range("J122").Value = 1777
$ws.Range("K122").Value = 6533.25
$ws.Range("L122").Value = 5331
$ws.Range("M122").Value = -4083.25
$ws.Range("N122").Value = -10231
$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2039.1538
$ws.Range("I132").Value = 1960.72
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5882.16
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3352.16
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3744.25
$ws.Range("I40").Value = 3497.1667
$ws.Range("J40").Value = 3991.3333
$ws.Range("K40").Value = 3497.1667
$ws.Range("L40").Value = 3991.3333
$ws.Range("M40").Value = -3361.1667
$ws.Range("N40").Value = -4263.3333
$ws.Range("H69").Value = 50777
$ws.Range("J69").Value = 50777
$ws.Range("L69").Value = 50777
$ws.Range("N69").Value = -52399
$ws.Range("H72").Value = 50777
$ws.Range("J72").Value = 50777
$ws.Range("L72").Value = 152331
$ws.Range("N72").Value = -160443
$ws.Range("H93").Value = 6372.9316
$ws.Range("I93").Value = 1671.7667
$ws.Range("J93").Value = 16446.857
$ws.Range("K93").Value = 1671.7667
$ws.Range("L93").Value = 16446.857
$ws.Range("M93").Value = -423.7666999999999
$ws.Range("N93").Value = -18942.857
$ws.Range("H132").Value = 8146.5244
$ws.Range("I132").Value = 7952.904
$ws.Range("J132").Value = 9265.223
$ws.Range("K132").Value = 23858.712
$ws.Range("L132").Value = 27795.669
$ws.Range("M132").Value = -21328.712
$ws.Range("N132").Value = -32855.669
$ws.Range("H136").Value = 3897.9148
$ws.Range("I136").Value = 4004.8809
$ws.Range("J136").Value = 2999.4
$ws.Range("K136").Value = 12014.6427
$ws.Range("L136").Value = 8998.200000000001
$ws.Range("M136").Value = -9464.6427
$ws.Range("N136").Value = -14098.2
$ws.Range("H137").Value = 85095
$ws.Range("J137").Value = 85095
$ws.Range("L137").Value = 85095
$ws.Range("N137").Value = -95295

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2131708
$ws.Range("J5").Value = 2131708
$ws.Range("L5").Value = 2131708
$ws.Range("N5").Value = -2131932
$ws.Range("H15").Value = 7487.6665
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H51").Value = 49999
$ws.Range("J51").Value = 49999
$ws.Range("L51").Value = 49999
$ws.Range("N51").Value = -51019
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()
$ws.Range("H64").Value = 65000
$ws.Range("J64").Value = 65000
$ws.Range("L64").Value = 65000
$ws.Range("N64").Value = -65496
$ws.Range("H67").Value = 65000
$ws.Range("J67").Value = 65000
$ws.Range("L67").Value = 65000
$ws.Range("N67").Value = -66716
$ws.Range("H70").Value = 40763.5
$ws.Range("J70").Value = 40777
$ws.Range("L70").Value = 40777
$ws.Range("N70").Value = -41407
$ws.Range("H73").Value = 40763.5
$ws.Range("J73").Value = 40777
$ws.Range("L73").Value = 40777
$ws.Range("N73").Value = -42961
$ws.Range("H100").Value = 933.4
$ws.Range("J100").Value = 1159
$ws.Range("L100").Value = 2318
$ws.Range("N100").Value = -3400
$ws.Range("H101").Value = 20867.334
$ws.Range("J101").Value = 20867.334
$ws.Range("L101").Value = 20867.334
$ws.Range("N101").Value = -27357.334
$ws.Range("H122").Value = 2875.5173
$ws.Range("I122").Value = 1538.5217
$ws.Range("J122").Value = 8000.6665
$ws.Range("K122").Value = 4615.5651
$ws.Range("L122").Value = 24001.9995
$ws.Range("M122").Value = -2165.5651
$ws.Range("N122").Value = -28901.9995
$ws.Range("H126").Value = 3691.125
$ws.Range("I126").Value = 4083.077
$ws.Range("K126").Value = 12249.231
$ws.Range("M126").Value = -9779.231
$ws.Range("H132").Value = 2346.1404
$ws.Range("I132").Value = 2455.0815
$ws.Range("J132").Value = 1678.875
$ws.Range("K132").Value = 7365.244499999999
$ws.Range("L132").Value = 5036.625
$ws.Range("M132").Value = -4835.244499999999
$ws.Range("N132").Value = -10096.625
$ws.Range("H136").Value = 5809.1304
$ws.Range("I136").Value = 3774.4666
$ws.Range("K136").Value = 11323.3998
$ws.Range("M136").Value = -8773.399800000001
